$d = $word.ActiveDocument

# Locate the paragraph containing "Page A propos" robustly (search by text,
# then resolve to the Paragraphs collection so we can grab the paragraph
# that follows it, i.e. "Footer").
$rng = $d.Content
$found = $rng.Find.Execute("Page A propos", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find paragraph 'Page A propos'"
}

$allParas = $d.Paragraphs
$idx = -1
for ($i = 1; $i -le $allParas.Count; $i++) {
    if ($allParas.Item($i).Range.Start -eq $rng.Start) {
        $idx = $i
        break
    }
}
if ($idx -eq -1) {
    throw "Could not resolve paragraph index for 'Page A propos'"
}

$pPageAPropos = $allParas.Item($idx)
$pFooter = $allParas.Item($idx + 1)

if ($pFooter.Range.Text.TrimEnd([char]13) -ne "Footer") {
    throw "Unexpected paragraph following 'Page A propos': $($pFooter.Range.Text)"
}

# Whole range covering both the "Page A propos" paragraph and the "Footer"
# paragraph (including their paragraph marks and the trailing _GoBack bookmark).
$full = $d.Range($pPageAPropos.Range.Start, $pFooter.Range.End)

$xml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
  <w:pPr>
    <w:pStyle w:val="Paragraphedeliste"/>
    <w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr>
    <w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:lang w:val="pt-PT"/></w:rPr>
  </w:pPr>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:lang w:val="pt-PT"/></w:rPr>
    <w:lastRenderedPageBreak/>
    <w:t>Page</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:lang w:val="pt-PT"/></w:rPr>
    <w:t xml:space="preserve"> A </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:lang w:val="pt-PT"/></w:rPr>
    <w:t>propos</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:lang w:val="pt-PT"/></w:rPr>
    <w:t xml:space="preserve"> + </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:lang w:val="pt-PT"/></w:rPr>
    <w:t>Forum</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:lang w:val="pt-PT"/></w:rPr>
    <w:t xml:space="preserve"> =&gt; </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:lang w:val="pt-PT"/></w:rPr>
    <w:t>F</w:t>
  </w:r>
  <w:r>
    <w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:lang w:val="pt-PT"/></w:rPr>
    <w:t>ooter</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="Paragraphedeliste"/>
    <w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr>
    <w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr>
  </w:pPr>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr>
    <w:t>Footer</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="Paragraphedeliste"/>
    <w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr>
    <w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr>
    <w:t>Affichage identifiants invalides</w:t>
  </w:r>
  <w:r>
    <w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr>
    <w:t xml:space="preserve"> (mode t&#233;l&#233;phone)</w:t>
  </w:r>
  <w:bookmarkStart w:id="0" w:name="_GoBack"/>
  <w:bookmarkEnd w:id="0"/>
</w:p>
</w:body></w:document>
</pkg:xmlData></pkg:part>
</pkg:package>
"@

[void]$full.InsertXML($xml)
